# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
#
# Semantic change (per the target diff): cell B11 on the "Rules" sheet
# changes from the text "R40" to the text "1" (still a plain text value,
# not a number), while keeping its existing style (s="23") untouched and
# adding the new literal "1" to the shared-strings table.
#
# Simply assigning Range.Value = "1" would let Excel auto-coerce the
# numeric-looking text into a real number (and pick up a different cell
# format along the way). To force it to stay text - exactly like the
# original "R40" was text - we stage the literal through a scratch cell
# using a formula that evaluates to the text string "1", copy that
# computed text value, and paste only the *value* into B11 so the
# destination's number format / style is left completely alone.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$helper = $ws.Range("Z1")
$helper.Formula = "=""1"""      # formula that evaluates to the text "1"
$helper.Copy()
$ws.Range("B11").PasteSpecial(-4163)   # xlPasteValues: value only, keep B11's own formatting
$helper.Clear()                 # remove the scratch cell and its formatting entirely
